$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds numeric-looking text (e.g. "157.66"); Excel would otherwise
# auto-convert these to real numbers on assignment. Force text format on just
# the cells we touch so they keep their original General/text representation.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '72.572.54'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').Value = '3.971.81'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '585.08'
$ws.Range('E5').Value = '  +8.39%  '
$ws.Range('D6').Value = '157.66'
$ws.Range('E6').Value = '  +6.39%  '
$ws.Range('D7').Value = '0.679'
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('D9').Value = '0.748'
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('D11').Value = '53.24'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = '0.0000318'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '10.86'
$ws.Range('E13').Value = '  +3.56%  '
$ws.Range('D14').Value = '4.604.90'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = '3.959.04'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('E16').Value = '  +9.40%  '
$ws.Range('D17').Value = '14.02'
$ws.Range('E17').Value = '  +1.09%  '
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '72.253.70'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('D21').Value = '433.21'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('D22').Value = '4.67'
$ws.Range('E22').Value = '  +10.45%  '
$ws.Range('D23').Value = '95.90'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = '3.42'
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '4.48'
$ws.Range('E25').Value = '  +24.46%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '14.31'
$ws.Range('E26').Value = '  +1.11%  '
$ws.Range('D27').Value = '11.10'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').Value = '10.70'
$ws.Range('E28').Value = '  +2.34%  '
$ws.Range('D29').Value = '5.92'
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('D30').Value = '36.41'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').Value = '7.83'
$ws.Range('E31').Value = '  +6.68%  '
$ws.Range('D32').Value = '13.57'
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('E33').Value = '  +1.72%  '
$ws.Range('D34').Value = '48.86'
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('D35').Value = '680.41'
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').Value = '68.98'
$ws.Range('E36').Value = '  +5.48%  '
$ws.Range('D37').Value = '0.435'
$ws.Range('E37').Value = '  +1.49%  '
$ws.Range('D38').Value = '0.0₃0856'
$ws.Range('E38').Value = '  +5.07%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').Value = '3.35'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').Value = '3.35'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.145'
$ws.Range('E42').Value = '  -2.19%  '
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').Value = '10.78'
$ws.Range('E44').Value = '  +11.66%  '
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('D47').Value = '2.65'
$ws.Range('E47').Value = '  -0.65%  '
$ws.Range('D48').Value = '3.36'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '3.02'
$ws.Range('E49').Value = '  +2.50%  '
$ws.Range('B50').Value = 'LidoDAOToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D50').Value = '3.42'
$ws.Range('E50').Value = '  +5.69%  '
$ws.Range('D51').Value = '2.16'
$ws.Range('E51').Value = '  +9.11%  '
